$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 291
$ws1.Range("F4").Value = 1207
$ws1.Range("F5").Value = 840
$ws1.Range("F6").Value = 867
$ws1.Range("F7").Value = 1582
$ws1.Range("F8").Value = 325
$ws1.Range("F9").Value = 1081
$ws1.Range("F14").Value = 545
$ws1.Range("F15").Value = 87
$ws1.Range("F17").Value = 16
$ws1.Range("F24").Value = 796

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 1051
$ws2.Range("F4").Value = 291
$ws2.Range("F6").Value = 209

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 275

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 275
$ws4.Range("F3").Value = 291
$ws4.Range("F4").Value = 1051
$ws4.Range("F6").Value = 1207
$ws4.Range("F7").Value = 840
$ws4.Range("F8").Value = 867
$ws4.Range("F9").Value = 1582
$ws4.Range("F10").Value = 325
$ws4.Range("F11").Value = 1081
$ws4.Range("F16").Value = 545
$ws4.Range("F17").Value = 87
$ws4.Range("F20").Value = 16
$ws4.Range("F21").Value = 291
$ws4.Range("F25").Value = 209
$ws4.Range("F26").Value = 209
$ws4.Range("F31").Value = 796
